$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56, shifting existing rows 56-59 down to 57-60
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly data point
$ws.Cells.Item(56, 1).Value = 11
$ws.Cells.Item(56, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(56, 3).Value = "Bíobío"
$ws.Cells.Item(56, 4).Value = 45127
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = 8
$ws.Cells.Item(56, 6).Value = 100114007
$ws.Cells.Item(56, 7).Value = "Jengibre"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 35
$ws.Cells.Item(56, 11).Value = 17000
$ws.Cells.Item(56, 12).Value = 18000
$ws.Cells.Item(56, 13).Value = 17429
$ws.Cells.Item(56, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(56, 15).Value = "Perú"
$ws.Cells.Item(56, 16).Value = 1341
$ws.Cells.Item(56, 17).Value = 13
$ws.Cells.Item(56, 18).Value = "Hortaliza"
